$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - headers
$ws.Range("A1").Value = "Technology"
$ws.Range("B1").Value = "Layout"
$ws.Range("C1").Value = "Lot"
$ws.Range("D1").Value = "Wafer"
$ws.Range("E1").Value = "Yield"

# Format the Yield column as text so percentage-looking strings stay literal text
$ws.Range("E2:E4").NumberFormat = "@"

# Row 2
$ws.Range("C2").Value = "W118892"
$ws.Range("D2").Value = "R114792-03"
$ws.Range("E2").Value = "99.97%"

# Row 3
$ws.Range("C3").Value = "GAL-LOT"
$ws.Range("D3").Value = "GAL-LOT-02"
$ws.Range("E3").Value = "99.97%"

# Row 4
$ws.Range("C4").Value = "GAL-LOT"
$ws.Range("D4").Value = "GAL-LOT-03"
$ws.Range("E4").Value = "99.97%"

# The NumberFormat = "@" trick above keeps "99.97%" as literal text instead of
# auto-converting to a percentage number, but it leaves a style applied to the
# cells. Clear the formatting again so the cells end up with the default style
# (matching the original workbook, which has no per-cell styles).
$ws.Range("E2:E4").ClearFormats()
